$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.80"
$ws.Range("E2").Value = "'2.48%"
$ws.Range("B2:E2").Style = "Normal"

$ws.Range("D3").Value = "'39.50"
$ws.Range("E3").Value = "'2.68%"
$ws.Range("B3:E3").Style = "Normal"

$ws.Range("D4").Value = "'5.131"
$ws.Range("E4").Value = "'0.96%"
$ws.Range("B4:E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08195"
$ws.Range("E5").Value = "'1.27%"
$ws.Range("B5:E5").Style = "Normal"

$ws.Range("D6").Value = "'1.959"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("B6:E6").Style = "Normal"

$ws.Range("D7").Value = "'8.230"
$ws.Range("E7").Value = "'3.51%"
$ws.Range("B7:E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9277"
$ws.Range("E8").Value = "'0.01%"
$ws.Range("B8:E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1418"
$ws.Range("E9").Value = "'-3.25%"
$ws.Range("B9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1974"
$ws.Range("E10").Value = "'2.11%"
$ws.Range("B10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09139"
$ws.Range("E11").Value = "'0.18%"
$ws.Range("B11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03510"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("B12:E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09820"
$ws.Range("E13").Value = "'0.13%"
$ws.Range("B13:E13").Style = "Normal"

$ws.Range("E14").Value = "'-0.74%"
$ws.Range("B14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.005981"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("B15:E15").Style = "Normal"

$ws.Range("D16").Value = "'3.651"
$ws.Range("E16").Value = "'-1.95%"
$ws.Range("B16:E16").Style = "Normal"

$ws.Range("D17").Value = "'4.241"
$ws.Range("E17").Value = "'0.83%"
$ws.Range("B17:E17").Style = "Normal"

$ws.Range("D18").Value = "'3.230"
$ws.Range("E18").Value = "'-5.44%"
$ws.Range("B18:E18").Style = "Normal"

$ws.Range("D19").Value = "'0.3460"
$ws.Range("E19").Value = "'-0.07%"
$ws.Range("B19:E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1304"
$ws.Range("E20").Value = "'-2.30%"
$ws.Range("B20:E20").Style = "Normal"

$ws.Range("D21").Value = "'4.821"
$ws.Range("E21").Value = "'0.32%"
$ws.Range("B21:E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2445"
$ws.Range("E22").Value = "'-0.45%"
$ws.Range("B22:E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04359"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("B23:E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.66%"
$ws.Range("B24:E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004781"
$ws.Range("E25").Value = "'-1.11%"
$ws.Range("B25:E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001297"
$ws.Range("E26").Value = "'-0.38%"
$ws.Range("B26:E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003992"
$ws.Range("E27").Value = "'-10.25%"
$ws.Range("B27:E27").Style = "Normal"

$ws.Range("D39").Value = "'0.02206"
$ws.Range("E39").Value = "'5.86%"
$ws.Range("B39:E39").Style = "Normal"

$ws.Range("D40").Value = "'0.05214"
$ws.Range("E40").Value = "'3.19%"
$ws.Range("B40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007531"
$ws.Range("E41").Value = "'0.62%"
$ws.Range("B41:E41").Style = "Normal"

$ws.Range("D42").Value = "'0.009804"
$ws.Range("E42").Value = "'-3.37%"
$ws.Range("B42:E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1375"
$ws.Range("E43").Value = "'1.50%"
$ws.Range("B43:E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002107"
$ws.Range("E44").Value = "'-1.71%"
$ws.Range("B44:E44").Style = "Normal"

$ws.Range("D45").Value = "'0.009813"
$ws.Range("E45").Value = "'6.58%"
$ws.Range("B45:E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006356"
$ws.Range("E46").Value = "'2.46%"
$ws.Range("B46:E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.42%"
$ws.Range("B47:E47").Style = "Normal"

$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001198"
$ws.Range("E48").Value = "'-25.26%"
$ws.Range("B48:E48").Style = "Normal"

$ws.Range("B49").Value = "'BOLO"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002759"
$ws.Range("E49").Value = "'-7.60%"
$ws.Range("B49:E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.42%"
$ws.Range("B50:E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.42%"
$ws.Range("B51:E51").Style = "Normal"
